$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Title paragraph: "Entertainment Law" -> "Law & Business For Artists"
#    (scope the Find to the first paragraph only so nothing else is touched)
# ---------------------------------------------------------------------
$titleRange = $d.Paragraphs(1).Range
$titleRange.Find.Execute("Entertainment Law", $true, $false, $false, $false, $false,
                          $true, 1, $false, "Law & Business For Artists", 2)

# ---------------------------------------------------------------------
# 2. Subtitle paragraph: "Lawyers For The Talent" -> "Art Biz Law"
#    (scope the Find to the second paragraph only; the Heading1 paragraph
#    further down the document has the same text and must stay untouched)
# ---------------------------------------------------------------------
$subtitleRange = $d.Paragraphs(2).Range
$subtitleRange.Find.Execute("Lawyers For The Talent", $true, $false, $false, $false, $false,
                             $true, 1, $false, "Art Biz Law", 2)

# ---------------------------------------------------------------------
# 3. Remove the stray "Heading6" paragraph style from the four byline /
#    court-attribution paragraphs so they fall back to the Normal style.
#    Located via their bookmarks so the edit is robust to any paragraph
#    re-numbering caused by the edits above.
# ---------------------------------------------------------------------
$bookmarkNames = @(
    "by-richard-dooling",
    "moral-rights-droit-morale-and-other-rights",
    "united-states-court-of-appeals-second-circuit-1952.",
    "united-states-court-of-appeals-second-circuit-1976"
)
foreach ($name in $bookmarkNames) {
    $bm = $d.Bookmarks($name)
    $para = $bm.Range.Paragraphs(1)
    $para.Range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 4. Clarify the Herbert Harris sentence with parentheses.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Herbert Harris a record producer", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Herbert Harris (a record producer)", 2)

Write-Host "Edits applied"
